# Optimize the display of the results
# Replace the old holdings table (rows 2-11) with the new holdings table (rows 2-9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range below the header so stale rows (10-11) disappear
$ws.Range("A2:D11").Clear()

# Header row stays the same conceptually (Ticker/Cash, Share count, Cost basis, TOTAL VALUE)
$ws.Range("A1").Value = "Ticker/Cash"
$ws.Range("B1").Value = "Share count"
$ws.Range("C1").Value = "Cost basis"
$ws.Range("D1").Value = "TOTAL VALUE"

# New holdings data
$ws.Range("A2").Value = "HWM"
$ws.Range("B2").Value = 630
$ws.Range("C2").Value = 172.2

$ws.Range("A3").Value = "PLTR"
$ws.Range("B3").Value = 770
$ws.Range("C3").Value = 139.6

$ws.Range("A4").Value = "NFLX"
$ws.Range("B4").Value = 89
$ws.Range("C4").Value = 1231.6

$ws.Range("A5").Value = "NVDA"
$ws.Range("B5").Value = 1884
$ws.Range("C5").Value = 137.5

$ws.Range("A6").Value = "MSFT"
$ws.Range("B6").Value = 448
$ws.Range("C6").Value = 457

$ws.Range("A7").Value = "KLAC"
$ws.Range("B7").Value = 124
$ws.Range("C7").Value = 774.2

$ws.Range("A8").Value = "APH"
$ws.Range("B8").Value = 1730
$ws.Range("C8").Value = 93.6

$ws.Range("A9").Value = "Cash"
$ws.Range("D9").Value = 7490

# Update the view: zoom in and move the active selection to C10
$ws.Activate()
$excel.ActiveWindow.Zoom = 333
$ws.Range("C10").Select()
